# Update worksheet values to reflect newly computed TPM-based NATMI metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 0.1635273785303319
$ws.Range("J2").Value = 0.1635273785303318
$ws.Range("M2").Value = 1.667434
$ws.Range("N2").Value = 5.002302
$ws.Range("O2").Value = 0.3223739883484499
$ws.Range("P2").Value = 0.32237398834845
$ws.Range("Q2").Value = 0.004154689716666667
$ws.Range("R2").Value = 0.03739220745
$ws.Range("S2").Value = 0.05271697322098977
$ws.Range("T2").Value = 0.05271697322098977
$ws.Range("I3").Value = 0.1635273785303319
$ws.Range("J3").Value = 0.1635273785303318
$ws.Range("O3").Value = 0.2193354457157105
$ws.Range("P3").Value = 0.2193354457157106
$ws.Range("S3").Value = 0.03586735045667205
$ws.Range("T3").Value = 0.03586735045667205
$ws.Range("I4").Value = 0.1635273785303319
$ws.Range("J4").Value = 0.1635273785303318
$ws.Range("M4").Value = 0.2055123333333333
$ws.Range("N4").Value = 0.616537
$ws.Range("O4").Value = 0.03973280534729575
$ws.Range("P4").Value = 0.03973280534729576
$ws.Range("Q4").Value = 0.0005120682305555555
$ws.Range("R4").Value = 0.004608614075
$ws.Range("S4").Value = 0.006497401500099226
$ws.Range("T4").Value = 0.006497401500099228
$ws.Range("I5").Value = 0.1635273785303319
$ws.Range("J5").Value = 0.1635273785303318
$ws.Range("M5").Value = 1.661741333333333
$ws.Range("N5").Value = 4.985224
$ws.Range("O5").Value = 0.3212733944672698
$ws.Range("P5").Value = 0.3212733944672699
$ws.Range("Q5").Value = 0.004140505488888889
$ws.Range("R5").Value = 0.0372645494
$ws.Range("S5").Value = 0.05253699598877386
$ws.Range("T5").Value = 0.05253699598877386
$ws.Range("I6").Value = 0.1635273785303319
$ws.Range("J6").Value = 0.1635273785303318
$ws.Range("M6").Value = 0.1178836666666667
$ws.Range("N6").Value = 0.353651
$ws.Range("O6").Value = 0.02279108365576842
$ws.Range("P6").Value = 0.02279108365576842
$ws.Range("Q6").Value = 0.0002937268027777778
$ws.Range("R6").Value = 0.002643541225
$ws.Range("S6").Value = 0.003726966164093301
$ws.Range("T6").Value = 0.003726966164093301
$ws.Range("I7").Value = 0.1635273785303319
$ws.Range("J7").Value = 0.1635273785303318
$ws.Range("M7").Value = 0.385306
$ws.Range("N7").Value = 1.155918
$ws.Range("O7").Value = 0.07449328246550557
$ws.Range("P7").Value = 0.0744932824655056
$ws.Range("Q7").Value = 0.0009600541166666667
$ws.Range("R7").Value = 0.008640487049999999
$ws.Range("S7").Value = 0.01218169119970366
$ws.Range("T7").Value = 0.01218169119970366
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.01274533333333334
$ws.Range("H8").Value = 0.03823600000000001
$ws.Range("I8").Value = 0.8364726214696682
$ws.Range("J8").Value = 0.8364726214696681
$ws.Range("M8").Value = 1.667434
$ws.Range("N8").Value = 5.002302
$ws.Range("O8").Value = 0.3223739883484499
$ws.Range("P8").Value = 0.32237398834845
$ws.Range("Q8").Value = 0.02125200214133334
$ws.Range("R8").Value = 0.191268019272
$ws.Range("S8").Value = 0.2696570151274602
$ws.Range("T8").Value = 0.2696570151274602
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.01274533333333334
$ws.Range("H9").Value = 0.03823600000000001
$ws.Range("I9").Value = 0.8364726214696682
$ws.Range("J9").Value = 0.8364726214696681
$ws.Range("O9").Value = 0.2193354457157105
$ws.Range("P9").Value = 0.2193354457157106
$ws.Range("Q9").Value = 0.01445934700222223
$ws.Range("R9").Value = 0.13013412302
$ws.Range("S9").Value = 0.1834680952590385
$ws.Range("T9").Value = 0.1834680952590385
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.01274533333333334
$ws.Range("H10").Value = 0.03823600000000001
$ws.Range("I10").Value = 0.8364726214696682
$ws.Range("J10").Value = 0.8364726214696681
$ws.Range("M10").Value = 0.2055123333333333
$ws.Range("N10").Value = 0.616537
$ws.Range("O10").Value = 0.03973280534729575
$ws.Range("P10").Value = 0.03973280534729576
$ws.Range("Q10").Value = 0.002619323192444445
$ws.Range("R10").Value = 0.023573908732
$ws.Range("S10").Value = 0.03323540384719652
$ws.Range("T10").Value = 0.03323540384719653
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.01274533333333334
$ws.Range("H11").Value = 0.03823600000000001
$ws.Range("I11").Value = 0.8364726214696682
$ws.Range("J11").Value = 0.8364726214696681
$ws.Range("M11").Value = 1.661741333333333
$ws.Range("N11").Value = 4.985224
$ws.Range("O11").Value = 0.3212733944672698
$ws.Range("P11").Value = 0.3212733944672699
$ws.Range("Q11").Value = 0.02117944720711111
$ws.Range("R11").Value = 0.190615024864
$ws.Range("S11").Value = 0.268736398478496
$ws.Range("T11").Value = 0.268736398478496
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.01274533333333334
$ws.Range("H12").Value = 0.03823600000000001
$ws.Range("I12").Value = 0.8364726214696682
$ws.Range("J12").Value = 0.8364726214696681
$ws.Range("M12").Value = 0.1178836666666667
$ws.Range("N12").Value = 0.353651
$ws.Range("O12").Value = 0.02279108365576842
$ws.Range("P12").Value = 0.02279108365576842
$ws.Range("Q12").Value = 0.001502466626222223
$ws.Range("R12").Value = 0.013522199636
$ws.Range("S12").Value = 0.01906411749167512
$ws.Range("T12").Value = 0.01906411749167512
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.01274533333333334
$ws.Range("H13").Value = 0.03823600000000001
$ws.Range("I13").Value = 0.8364726214696682
$ws.Range("J13").Value = 0.8364726214696681
$ws.Range("M13").Value = 0.385306
$ws.Range("N13").Value = 1.155918
$ws.Range("O13").Value = 0.07449328246550557
$ws.Range("P13").Value = 0.0744932824655056
$ws.Range("Q13").Value = 0.004910853405333334
$ws.Range("R13").Value = 0.044197680648
$ws.Range("S13").Value = 0.06231159126580191
$ws.Range("T13").Value = 0.06231159126580193
